$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "datetime" column (column A) was removed from the sheet.
# Select the whole column and delete it so every remaining column
# shifts one place to the left (B->A, C->B, ... K->J).
$ws.Columns("A:A").Select()
$ws.Columns("A:A").Delete()
